$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 301
$ws.Range("F3").Value = 1031
$ws.Range("F4").Value = 1235
$ws.Range("F5").Value = 1098
$ws.Range("F6").Value = 3267
$ws.Range("F9").Value = 1156
$ws.Range("F10").Value = 717
$ws.Range("F13").Value = 45
$ws.Range("F14").Value = 110
$ws.Range("F15").Value = 640
$ws.Range("F16").Value = 1533
$ws.Range("F17").Value = 1533
$ws.Range("F18").Value = 18
$ws.Range("F20").Value = 28
$ws.Range("F21").Value = 597
$ws.Range("F22").Value = 327
$ws.Range("F23").Value = 542
$ws.Range("F24").Value = 564
$ws.Range("F25").Value = 48839
$ws.Range("F26").Value = 48844
$ws.Range("F27").Value = 708
$ws.Range("F28").Value = 633
$ws.Range("F29").Value = 32497
$ws.Range("F30").Value = 32497
$ws.Range("F31").Value = 435
$ws.Range("F32").Value = 9
$ws.Range("F33").Value = 1
$ws.Range("F35").Value = 3
$ws.Range("F36").Value = 919
$ws.Range("F37").Value = 226
$ws.Range("F38").Value = 146
$ws.Range("F39").Value = 491
$ws.Range("F40").Value = 1147
$ws.Range("F41").Value = 5304
$ws.Range("F42").Value = 696
$ws.Range("F43").Value = 414
$ws.Range("F46").Value = 318

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 15
$ws.Range("F13").Value = 1789
$ws.Range("F14").Value = 8
$ws.Range("F15").Value = 819
$ws.Range("F17").Value = 61
$ws.Range("F18").Value = 392
$ws.Range("F20").Value = 58
$ws.Range("F25").Value = 759
$ws.Range("F28").Value = 23
$ws.Range("F30").Value = 12
$ws.Range("F35").Value = 1196
$ws.Range("F40").Value = 1
$ws.Range("F45").Value = 801

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 691
$ws.Range("F5").Value = 528
$ws.Range("F6").Value = 524

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 691
$ws.Range("F3").Value = 301
$ws.Range("F5").Value = 1032
$ws.Range("F6").Value = 1235
$ws.Range("F8").Value = 1098
$ws.Range("F9").Value = 3268
$ws.Range("F12").Value = 1156
$ws.Range("F13").Value = 717
$ws.Range("F14").Value = 524
$ws.Range("F17").Value = 1789
$ws.Range("F19").Value = 1533
$ws.Range("F20").Value = 1533
$ws.Range("F21").Value = 18
$ws.Range("F24").Value = 61
$ws.Range("F25").Value = 28
$ws.Range("F26").Value = 597
$ws.Range("F27").Value = 392
$ws.Range("F28").Value = 327
$ws.Range("F29").Value = 564
$ws.Range("F31").Value = 58
$ws.Range("F32").Value = 48854
$ws.Range("F34").Value = 32497
$ws.Range("F35").Value = 9
$ws.Range("F37").Value = 919
$ws.Range("F39").Value = 226
$ws.Range("F40").Value = 146
$ws.Range("F41").Value = 23
$ws.Range("F42").Value = 491
$ws.Range("F43").Value = 1147
$ws.Range("F44").Value = 5304
$ws.Range("F45").Value = 696
$ws.Range("F47").Value = 414
$ws.Range("F50").Value = 318
